# ConfigurationInputsTrout.xlsx update - "Update Trout ConfigInputs based on optimization"
#
# This script:
#  1. Inserts two new "Burial Factor" rows (Right/Left) into the SWGW section
#  2. Inserts two new "POC_lc" rows (Right/Left) into the Leach/Resp section
#  3. Inserts two new "RespParam" rows (DOCR/DOCL) into the Leach/Resp section
#  4. Updates the R_auto value from optimization
#  5. Minor cosmetic tweaks: column widths, selection, page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert BurialFactor_R / BurialFactor_L rows before "#SWGW Parameters" (row 11) ---
$ws.Rows.Item(11).Resize(2).Insert()
$ws.Range("A11").Value2 = "BurialFactor_R"
$ws.Range("B11").Value2 = 0.0014696080000000001
$ws.Range("C11").Value2 = "1/days"

$ws.Range("A12").Value2 = "BurialFactor_L"
$ws.Range("B12").Value2 = 0.00083764389999999998
$ws.Range("C12").Value2 = "1/days"

# --- 2. Insert POC_lcR / POC_lcL rows after "POC_lc" (now row 23) ---
$ws.Rows.Item(24).Resize(2).Insert()
$ws.Range("A24").Value2 = "POC_lcR"
$ws.Range("B24").Value2 = 0.17578609049999999
$ws.Range("C24").Value2 = "1/days"

$ws.Range("A25").Value2 = "POC_lcL"
$ws.Range("B25").Value2 = 0.0095980883999999995
$ws.Range("C25").Value2 = "1/days"

# --- 3. Insert DOCR_RespParam / DOCL_RespParam rows after "RespParam" (now row 26) ---
$ws.Rows.Item(27).Resize(2).Insert()
$ws.Range("A27").Value2 = "DOCR_RespParam"
$ws.Range("B27").Value2 = 0.0009963455
$ws.Range("C27").Value2 = "unitless"

$ws.Range("A28").Value2 = "DOCL_RespParam"
$ws.Range("B28").Value2 = 0.12293787389999999
$ws.Range("C28").Value2 = "unitless"

# --- 4. Update R_auto value (now row 29) with the optimized result ---
$ws.Range("B29").Value2 = 0.78708956070000002

# --- 5. Cosmetic adjustments ---
# Column widths (closest achievable snap given engine's width granularity)
$ws.Columns.Item(1).ColumnWidth = 25.0416666
$ws.Columns.Item(3).ColumnWidth = 12.7916666

# Selection / scroll position
$ws.Range("B25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10

# Page orientation
$ws.PageSetup.Orientation = 1

Write-Host "Done applying ConfigInputs optimization update."
